$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '52.196.27'
Set-TextValue "E2" '  -0.23%  '
Set-TextValue "D3" '2.839.95'
Set-TextValue "E3" '  +1.70%  '
Set-TextValue "E4" '  +0.04%  '
Set-TextValue "D5" '360.90'
Set-TextValue "E5" '  +5.89%  '
Set-TextValue "D6" '113.26'
Set-TextValue "E6" '  -3.00%  '
Set-TextValue "D7" '0.576'
Set-TextValue "E7" '  +4.06%  '
Set-TextValue "E8" '  -0.02%  '
Set-TextValue "D9" '0.610'
Set-TextValue "E9" '  +4.96%  '
Set-TextValue "D10" '41.59'
Set-TextValue "E10" '  -1.12%  '
Set-TextValue "D11" '0.0864'
Set-TextValue "E11" '  -0.43%  '
Set-TextValue "E12" '  +1.18%  '
Set-TextValue "D13" '20.04'
Set-TextValue "E13" '  -0.28%  '
Set-TextValue "D14" '7.81'
Set-TextValue "E14" '  +2.25%  '
Set-TextValue "D15" '3.288.35'
Set-TextValue "E15" '  +1.69%  '
Set-TextValue "D16" '2.820.11'
Set-TextValue "E16" '  +0.88%  '
Set-TextValue "E17" '  +2.75%  '
Set-TextValue "D18" '52.121.92'
Set-TextValue "E18" '  -0.01%  '
Set-TextValue "D19" '7.61'
Set-TextValue "E19" '  +9.21%  '
Set-TextValue "E20" '  -1.84%  '
Set-TextValue "D21" '13.54'
Set-TextValue "E21" '  +1.43%  '
Set-TextValue "D22" '0.0₃0996'
Set-TextValue "D23" '70.43'
Set-TextValue "E23" '  +0.20%  '
Set-TextValue "D24" '268.04'
Set-TextValue "E24" '  -3.88%  '
Set-TextValue "E25" '  +2.51%  '
Set-TextValue "D26" '27.16'
Set-TextValue "E26" '  +1.02%  '
Set-TextValue "E27" '  -0.03%  '
Set-TextValue "D28" '10.42'
Set-TextValue "E28" '  +1.52%  '
Set-TextValue "E29" '  +1.17%  '
Set-TextValue "D30" '54.40'
Set-TextValue "E30" '  +7.89%  '
Set-TextValue "D31" '0.0487'
Set-TextValue "E31" '  +28.77%  '
Set-TextValue "E32" '  -1.31%  '
Set-TextValue "D33" '34.71'
Set-TextValue "E33" '  -0.42%  '
Set-TextValue "E34" '  +1.88%  '
Set-TextValue "D35" '5.53'
Set-TextValue "E35" '  +11.21%  '
Set-TextValue "E36" '  +2.01%  '
Set-TextValue "E37" '  +0.09%  '
Set-TextValue "E38" '  +1.00%  '
Set-TextValue "D39" '2.07'
Set-TextValue "E39" '  -2.29%  '
Set-TextValue "E40" '  -3.03%  '
Set-TextValue "D41" '23.94'
Set-TextValue "E41" '  +1.85%  '
Set-TextValue "E42" '  +1.24%  '
Set-TextValue "D43" '128.01'
Set-TextValue "E43" '  +2.59%  '
Set-TextValue "D44" '2.55'
Set-TextValue "E44" '  -7.33%  '
Set-TextValue "E45" '  -1.73%  '
Set-TextValue "D46" '3.42'
Set-TextValue "E46" '  +2.67%  '
Set-TextValue "D47" '2.111.30'
Set-TextValue "E47" '  +0.51%  '
Set-TextValue "D49" '1.00'
Set-TextValue "E49" '  +11.35%  '
Set-TextValue "E50" '  +5.24%  '
Set-TextValue "D51" '62.08'
Set-TextValue "E51" '  +3.12%  '
